$wb = $excel.ActiveWorkbook

# 1. Rename header cells on the existing sheets
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 2. Add a new worksheet that will become "PO Forecast"
$wsForecast = $wb.Worksheets.Add()
$wsForecast.Name = "PO Forecast"

# 3. Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# 4. Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45256.99999999999, 215, -127.2819616240443, 540.0302752828569),
    @(45263.99999999999, 222, -112.0177652381434, 566.7113705355395),
    @(45270.99999999999, 230, -75.09374613983702, 550.7573646980525),
    @(45277.99999999999, 238, -93.77374149752971, 561.1118324572782),
    @(45298.99999999999, 261, -76.11816055194485, 594.9531921569156),
    @(45305.99999999999, 268, -31.48768033369026, 576.6744575358828),
    @(45312.99999999999, 276, -66.01809145487717, 595.0599369930266),
    @(45326.99999999999, 291, -31.76160686019109, 632.4098518110344),
    @(45396.99999999999, 367, 39.36056411525289, 706.8862167029373),
    @(45410.99999999999, 383, 64.84621879486514, 705.9522490927834),
    @(45417.99999999999, 390, 59.25552596252086, 697.6275490908578),
    @(45424.99999999999, 398, 45.66900240934702, 716.5418546319248),
    @(45431.99999999999, 405, 70.34272704212167, 734.6804977435692),
    @(45438.99999999999, 413, 83.58588551594433, 755.2083048776033),
    @(45445.99999999999, 421, 64.86048560201893, 750.6031937591596),
    @(45452.99999999999, 428, 110.3274791492746, 751.6895232598628),
    @(45459.99999999999, 436, 93.8722924717601, 750.2607463193566),
    @(45466.99999999999, 444, 95.19092865434224, 763.0798670089423),
    @(45480.99999999999, 459, 142.5549600108583, 803.7473958788117),
    @(45487.99999999999, 466, 149.3898380470835, 803.9304711093868),
    @(45515.99999999999, 497, 176.0140198746305, 828.9865970976492),
    @(45522.99999999999, 505, 210.8191746835069, 826.2862035952419),
    @(45536.99999999999, 520, 208.6225336380417, 838.1100381326746),
    @(45543.99999999999, 527, 228.526751580412, 879.3620737638367),
    @(45550.99999999999, 535, 220.7387808438695, 852.6832003779086),
    @(45557.99999999999, 543, 214.1907005984134, 864.8464439357299),
    @(45564.99999999999, 550, 252.2000600420939, 870.7951470130828),
    @(45571.99999999999, 558, 213.084745060066, 888.9694937064312),
    @(45578.99999999999, 566, 259.4297224308351, 899.5986644404902),
    @(45585.99999999999, 573, 243.4601812033004, 904.9318520572649),
    @(45592.99999999999, 581, 270.3447125118425, 900.8729465473905)

)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$wsForecast.Range("A2:A32").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# 5. Move the new sheet to be last (after "Monthly Trend").
#    Re-fetch the worksheet handles by name right before moving: adding a
#    sheet shifts positional indices, and stale handles captured earlier
#    can resolve to the wrong sheet after the collection changes.
$wsForecast = $wb.Worksheets.Item("PO Forecast")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsForecast.Move($null, $wsMonthly)
